$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row to append (row 25), mirroring the existing data rows
$row = 25

$ws.Cells.Item($row, 1).Value = 5
$ws.Cells.Item($row, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value = "Maule"

# Date column - copy style/number format from the row above and set the value
$ws.Cells.Item($row, 4).Value = 44628
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 7
$ws.Cells.Item($row, 6).Value = 100112043
$ws.Cells.Item($row, 7).Value = "Pepino dulce"
$ws.Cells.Item($row, 8).Value = "Cultivar IV Región"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 15000
$ws.Cells.Item($row, 12).Value = 15000
$ws.Cells.Item($row, 13).Value = 15000
$ws.Cells.Item($row, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 833
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"
